$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "ReasonToReject" column (J) previously had no value for the rejected
# rows (3-9). Populate it with "Nil" for each of those rows.
$ws.Range("J3:J9").Value = "Nil"

# Leave the selection where the user last clicked after typing the values.
$ws.Range("J10").Select()
